$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Feedback ID"
$ws.Range("B1").Value = "User ID"
$ws.Range("C1").Value = "Comments"
$ws.Range("D1").Value = "Rating"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 101
$ws.Range("C2").Value = "loved the design!"
$ws.Range("D2").Value = 5

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 102
$ws.Range("C3").Value = "very good"
$ws.Range("D3").Value = 4

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 103
$ws.Range("C4").Value = "Didn’t like it much"
$ws.Range("D4").Value = 1

# Remove old rows 5 and 6 (now not part of the table)
$ws.Range("A5:D6").Delete()

$ws.Range("G5").Select()
